$d = $word.ActiveDocument

# wdReplace constants
$wdReplaceAll = 2
$wdFindWrap = 1  # wdFindContinue, effectively no wrap beyond the scoped range

# 1. Title/heading text "Play Gifts of Fortune Slot Game for Free 2021" -> "Play Gifts of Fortune Slot Game for Free"
#    This exact string occurs twice (the Heading1 title and the bold call-to-action near the
#    end) and both need the same replacement, so a document-wide replace-all is safe.
$d.Content.Find.Execute("Play Gifts of Fortune Slot Game for Free 2021", $true, $false, $false, $false, $false, $true, 1, $false, "Play Gifts of Fortune Slot Game for Free", $wdReplaceAll)

# 2. "What we like" bullet list - the phrase "117,649 ways to win" appears several times in the
#    document (inside longer sentences too), so scope the Find to the exact bullet paragraph to
#    avoid touching the other occurrences.
$p = $d.Paragraphs.Item(52)
$p.Range.Find.Execute("117,649 ways to win", $true, $false, $false, $false, $false, $true, 1, $false, "Inspired by traditional Chinese beliefs about luck and fortune", $wdReplaceAll)

# 3. "Original bonus features" -> "Up to 117,649 ways to win"
$p = $d.Paragraphs.Item(53)
$p.Range.Find.Execute("Original bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "Up to 117,649 ways to win", $wdReplaceAll)

# 4. "Free spins and lucky prizes" -> "Exciting bonus features including free spins and multiplier prizes"
$p = $d.Paragraphs.Item(55)
$p.Range.Find.Execute("Free spins and lucky prizes", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting bonus features including free spins and multiplier prizes", $wdReplaceAll)

# 5. "What we don't like" bullet list
# "No progressive jackpot" -> "Limited betting options"
$p = $d.Paragraphs.Item(57)
$p.Range.Find.Execute("No progressive jackpot", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting options", $wdReplaceAll)

# 6. "No gamble feature" -> "Bonus features may be difficult to trigger"
$p = $d.Paragraphs.Item(58)
$p.Range.Find.Execute("No gamble feature", $true, $false, $false, $false, $false, $true, 1, $false, "Bonus features may be difficult to trigger", $wdReplaceAll)

# 7. Closing italic blurb
$p = $d.Paragraphs.Item(60)
$p.Range.Find.Execute("Experience Chinese inspired slot gaming with Gifts of Fortune. Play free and experience the unique bonuses, mobile optimization and up to 117,649 ways to win.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the Gifts of Fortune slot game with up to 117,649 ways to win. Play for free and experience exciting bonus features.", $wdReplaceAll)
